$wb = $excel.ActiveWorkbook

# --- Sheet "C": column B formulas now derive the "init" label directly
#     from column A of the same sheet, instead of pulling the raw value
#     from R_input. ---
$wsC = $wb.Worksheets.Item("C")
for ($r = 1; $r -le 60; $r++) {
    $wsC.Cells.Item($r, 2).Formula = "=IF('C'!`$A$r=`"`",`"`",'C'!`$A$r&`"init`")"
}

# --- Sheet "C0": column B formulas revert to pulling the raw value from
#     R_input directly (instead of mirroring sheet "C"). ---
$wsC0 = $wb.Worksheets.Item("C0")
for ($r = 1; $r -le 60; $r++) {
    $wsC0.Cells.Item($r, 2).Formula = "=IF(R_input!`$C$r=0,`"`",R_input!`$C$r)"
}

# --- Restore the last-used selection on each affected sheet, then
#     return focus to the sheet that was active before these edits. ---
$wsC.Range("D6").Select()
$wsC0.Range("F7").Select()
$wb.Worksheets.Item("R_input").Activate()
